$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "AFTNPNN8Q880"
$ws.Range("A3").Value = "Z60XTEFFEGXY"
$ws.Range("A4").Value = "XM4AM2PNJY0C"
$ws.Range("A8").Select() | Out-Null
